$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: update candidate number (text-formatted) and name
$ws.Range("A2").Value = "49.323"
$ws.Range("D2").Value = "Bui Quang Tuan "

# Row 3: update candidate number, name, class, and chuyen khoa column
$ws.Range("A3").Value = 49.323999999999998
$ws.Range("D3").Value = "Nguyen Anh Tuan "
$ws.Range("E3").Value = "B12D49"
$ws.Range("F3").Value = 3

# Row 4 (new): add a new student record
$ws.Range("A4").Value = 49.325000000000003
$ws.Range("B4").Value = 2019
$ws.Range("C4").Value = 2020
$ws.Range("D4").Value = "Nguyen Thi Lam Vien"
$ws.Range("E4").Value = "B13D49"
$ws.Range("F4").Value = 4

# Update selection to reflect the new active cell after data entry
$null = $ws.Range("G4").Select()
